$d = $word.ActiveDocument

# --- Change 1: add "17/5/2023" after "Date Closed: " ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Date Closed:*") {
        $r = $p.Range
        $insertPos = $r.End - 1
        $insertRange = $d.Range($insertPos, $insertPos)
        $insertRange.Text = "17/5/2023"

        $newRange = $d.Range($insertPos, $insertPos + 9)
        # Force the new text into its own run (distinct from the preceding
        # space run) by nudging the size before setting the final value.
        $newRange.Font.Size = 99
        $newRange.Font.Bold = $true
        $newRange.Font.Size = 20
        break
    }
}

# --- Change 2: mark the "Screenshot after fixing" picture run as NoProof ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = $true
    }
}
